$wb = $excel.ActiveWorkbook

# Update the Schema sheet header strings (date changed from 15:30:10 to 23:57:33)
$ws = $wb.Worksheets.Item("!!_Schema")
$ws.Unprotect()
$ws.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 23:57:33'"
$ws.Range("A2").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-03-09 23:57:33' objTablesVersion='0.0.8'"
$ws.Protect()

# Update the Compound data sheet header string: reorder attributes (tableFormat earlier) + new date
$ws = $wb.Worksheets.Item("!!Compound")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' id='Compound' name='Compound' description='Compound' date='2020-03-09 23:57:33' objTablesVersion='0.0.8'"
$ws.Protect()

# Update the Model data sheet header string: reorder attributes (tableFormat earlier) + new date
$ws = $wb.Worksheets.Item("!!Model")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables type='Data' tableFormat='column' id='Model' name='Model' description='Model' date='2020-03-09 23:57:33' objTablesVersion='0.0.8'"
$ws.Protect()

# Update the Reaction data sheet header string: reorder attributes (tableFormat earlier) + new date
$ws = $wb.Worksheets.Item("!!Reaction")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' id='Reaction' name='Reaction' description='Reaction' date='2020-03-09 23:57:33' objTablesVersion='0.0.8'"
$ws.Protect()
